$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPUresources-2021-2022-2023")

# Add the summary table formulas in column F
$ws.Range("F16").Formula = "=F17"
$ws.Range("F15").Formula = "=F20-F16-F17"

# Leave the selection on F16, matching the final saved state
$ws.Range("F16").Select()
